$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.737.90'
$ws.Range("E2").Value = '  +1.24%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.875.80'
$ws.Range("E3").Value = '  +0.61%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.53'
$ws.Range("E5").Value = '  +0.73%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.72'
$ws.Range("E6").Value = '  +3.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.876.56'
$ws.Range("E7").Value = '  +0.66%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  +1.02%  '

$ws.Range("E10").Value = '  +2.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.51'
$ws.Range("E11").Value = '  +3.30%  '

$ws.Range("E12").Value = '  +1.38%  '

$ws.Range("E13").Value = '  +15.40%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.20'
$ws.Range("E14").Value = '  +1.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.518.75'
$ws.Range("E15").Value = '  +0.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.870.72'
$ws.Range("E16").Value = '  +1.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.739.55'
$ws.Range("E17").Value = '  +1.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.36'
$ws.Range("E18").Value = '  +1.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.43'
$ws.Range("E19").Value = '  +0.35%  '

$ws.Range("E20").Value = '  +0.66%  '

$ws.Range("E21").Value = '  +2.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '472.69'
$ws.Range("E22").Value = '  +1.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.734'
$ws.Range("E23").Value = '  +0.96%  '

$ws.Range("E24").Value = '  +1.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.01'
$ws.Range("E25").Value = '  +1.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.28'
$ws.Range("E26").Value = '  +2.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.28'
$ws.Range("E27").Value = '  +0.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.51'
$ws.Range("E28").Value = '  +5.11%  '

$ws.Range("E30").Value = '  +0.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.020.69'
$ws.Range("E31").Value = '  +0.75%  '

$ws.Range("E32").Value = '  +1.60%  '

$ws.Range("E33").Value = '  +0.85%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.41'
$ws.Range("E34").Value = '  +1.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.45'
$ws.Range("E35").Value = '  +0.83%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.834.84'
$ws.Range("E36").Value = '  +0.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.95'
$ws.Range("E37").Value = '  +21.76%  '

$ws.Range("E38").Value = '  +1.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.02'
$ws.Range("E39").Value = '  +1.96%  '

$ws.Range("E40").Value = '  +0.67%  '

$ws.Range("E41").Value = '  +0.72%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("E43").Value = '  +2.71%  '

$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("E47").Value = '  +3.14%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '420.39'
$ws.Range("E48").Value = '  -1.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '46.70'
$ws.Range("E49").Value = '  -0.94%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '142.70'
$ws.Range("E50").Value = '  -0.61%  '

$ws.Range("E51").Value = '  +1.59%  '

# Row 44 and 45 swap (Stacks <-> FLOKI) with updated values
$ws.Range("B44").Value = 'FLOKI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.000303'
$ws.Range("E44").Value = '  +14.11%  '

$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.01'
$ws.Range("E45").Value = '  +1.19%  '
